$wb = $excel.ActiveWorkbook

# Shared string header change: "Tm" -> "Team" on all sheets (cell D1)
foreach ($sheetName in @("per_game", "per_minute", "per_poss", "advanced")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("D1").Value = "Team"
}


# --- per_game ---
$ws = $wb.Worksheets.Item("per_game")
$ws.Range("G36").Value = 48
$ws.Range("I36").Value = 26.5
$ws.Range("L36").Value = 0.383
$ws.Range("O36").Value = 0.348
$ws.Range("Q36").Value = 1.1
$ws.Range("S36").Value = 0.5
$ws.Range("X36").Value = 3
$ws.Range("Z36").Value = 1.2
$ws.Range("AD36").Value = 2.5
$ws.Range("G38").Value = 16
$ws.Range("I38").Value = 19.6
$ws.Range("L38").Value = 0.444
$ws.Range("M38").Value = 0.8
$ws.Range("O38").Value = 0.462
$ws.Range("Q38").Value = 0.6
$ws.Range("S38").Value = 0.611
$ws.Range("Y38").Value = 2.9
$ws.Range("Z38").Value = 0.8
$ws.Range("AE38").Value = 2.8
$ws.Range("G43").Value = 16
$ws.Range("I43").Value = 19.6
$ws.Range("L43").Value = 0.444
$ws.Range("M43").Value = 0.8
$ws.Range("O43").Value = 0.462
$ws.Range("Q43").Value = 0.6
$ws.Range("S43").Value = 0.611
$ws.Range("Y43").Value = 2.9
$ws.Range("Z43").Value = 0.8
$ws.Range("AE43").Value = 2.8
$ws.Range("G50").Value = 702

# --- per_minute ---
$ws = $wb.Worksheets.Item("per_minute")
$ws.Range("G26").Value = 48
$ws.Range("I26").Value = 1273
$ws.Range("L26").Value = 0.383
$ws.Range("O26").Value = 0.348
$ws.Range("AC26").Value = 3.5
$ws.Range("G28").Value = 16
$ws.Range("I28").Value = 314
$ws.Range("L28").Value = 0.444
$ws.Range("M28").Value = 1.4
$ws.Range("N28").Value = 3
$ws.Range("O28").Value = 0.462
$ws.Range("Q28").Value = 1.1
$ws.Range("V28").Value = 1.9
$ws.Range("W28").Value = 3.4
$ws.Range("X28").Value = 5.4
$ws.Range("Y28").Value = 1.4
$ws.Range("Z28").Value = 1
$ws.Range("AB28").Value = 0.8
$ws.Range("AC28").Value = 3.6
$ws.Range("G33").Value = 16
$ws.Range("I33").Value = 314
$ws.Range("L33").Value = 0.444
$ws.Range("M33").Value = 1.4
$ws.Range("N33").Value = 3
$ws.Range("O33").Value = 0.462
$ws.Range("Q33").Value = 1.1
$ws.Range("V33").Value = 1.9
$ws.Range("W33").Value = 3.4
$ws.Range("X33").Value = 5.4
$ws.Range("Y33").Value = 1.4
$ws.Range("Z33").Value = 1
$ws.Range("AB33").Value = 0.8
$ws.Range("AC33").Value = 3.6
$ws.Range("G40").Value = 702
$ws.Range("I40").Value = 20465
$ws.Range("Y40").Value = 1.8
$ws.Range("Y42").Value = -0.3

# --- per_poss ---
$ws = $wb.Worksheets.Item("per_poss")
$ws.Range("G26").Value = 48
$ws.Range("I26").Value = 1273
$ws.Range("L26").Value = 0.383
$ws.Range("O26").Value = 0.348
$ws.Range("Q26").Value = 2
$ws.Range("V26").Value = 1.9
$ws.Range("X26").Value = 7.2
$ws.Range("Y26").Value = 2.1
$ws.Range("AA26").Value = 0.7
$ws.Range("N27").Value = 4.2
$ws.Range("G28").Value = 16
$ws.Range("I28").Value = 314
$ws.Range("L28").Value = 0.444
$ws.Range("N28").Value = 3.9
$ws.Range("O28").Value = 0.462
$ws.Range("Q28").Value = 1.5
$ws.Range("T28").Value = 0.1
$ws.Range("V28").Value = 2.5
$ws.Range("W28").Value = 4.5
$ws.Range("X28").Value = 7
$ws.Range("Y28").Value = 1.8
$ws.Range("Z28").Value = 1.3
$ws.Range("AD28").Value = 6.6
$ws.Range("AF28").Value = 125
$ws.Range("G33").Value = 16
$ws.Range("I33").Value = 314
$ws.Range("L33").Value = 0.444
$ws.Range("N33").Value = 3.9
$ws.Range("O33").Value = 0.462
$ws.Range("Q33").Value = 1.5
$ws.Range("T33").Value = 0.1
$ws.Range("V33").Value = 2.5
$ws.Range("W33").Value = 4.5
$ws.Range("X33").Value = 7
$ws.Range("Y33").Value = 1.8
$ws.Range("Z33").Value = 1.3
$ws.Range("AD33").Value = 6.6
$ws.Range("AF33").Value = 125
$ws.Range("G40").Value = 702
$ws.Range("I40").Value = 20465

# --- advanced ---
$ws = $wb.Worksheets.Item("advanced")
$ws.Range("G26").Value = 48
$ws.Range("H26").Value = 1273
$ws.Range("J26").Value = 0.521
$ws.Range("K26").Value = 0.671
$ws.Range("L26").Value = 0.144
$ws.Range("N26").Value = 12.4
$ws.Range("O26").Value = 8.1
$ws.Range("P26").Value = 5.7
$ws.Range("S26").Value = 18.4
$ws.Range("T26").Value = 7.2
$ws.Range("V26").Value = 0.1
$ws.Range("Y26").Value = 0.036
$ws.Range("P27").Value = 6.1
$ws.Range("W27").Value = 0.6
$ws.Range("Y27").Value = 0.023
$ws.Range("G28").Value = 16
$ws.Range("H28").Value = 314
$ws.Range("I28").Value = 6.4
$ws.Range("J28").Value = 0.604
$ws.Range("K28").Value = 0.722
$ws.Range("L28").Value = 0.028
$ws.Range("N28").Value = 9.699999999999999
$ws.Range("O28").Value = 7.9
$ws.Range("P28").Value = 4.4
$ws.Range("Q28").Value = 1.3
$ws.Range("S28").Value = 16.1
$ws.Range("T28").Value = 5.8
$ws.Range("V28").Value = 0.3
$ws.Range("AA28").Value = -2.5
$ws.Range("AB28").Value = 0.2
$ws.Range("G33").Value = 16
$ws.Range("H33").Value = 314
$ws.Range("I33").Value = 6.4
$ws.Range("J33").Value = 0.604
$ws.Range("K33").Value = 0.722
$ws.Range("L33").Value = 0.028
$ws.Range("N33").Value = 9.699999999999999
$ws.Range("O33").Value = 7.9
$ws.Range("P33").Value = 4.4
$ws.Range("Q33").Value = 1.3
$ws.Range("S33").Value = 16.1
$ws.Range("T33").Value = 5.8
$ws.Range("V33").Value = 0.3
$ws.Range("AA33").Value = -2.5
$ws.Range("AB33").Value = 0.2
$ws.Range("G40").Value = 702
$ws.Range("H40").Value = 20465
$ws.Range("I40").Value = 10.2
$ws.Range("K40").Value = 0.471
$ws.Range("L40").Value = 0.202
$ws.Range("S40").Value = 13
$ws.Range("I42").Value = 0.8000000000000007
$ws.Range("K42").Value = 0.2110000000000001
$ws.Range("L42").Value = -0.05500000000000002
$ws.Range("S42").Value = -2.4
